$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): labels for the OAuth/Okta register+login data
$ws.Range("A1").Value = "State"
$ws.Range("B1").Value = "Client ID"
$ws.Range("C1").Value = "Client secret"
$ws.Range("D1").Value = "Login"
$ws.Range("E1").Value = "Password"

# Data row (row 2): corresponding values
$ws.Range("A2").Value = "KKAITM7eldKlwLGn01qJ"
$ws.Range("B2").Value = "0oamk0pm9fQx125R10h7"
$ws.Range("C2").Value = "Vm2-JINkX0t3GYjrQyrbDtNNRhYbur5zA06dcZ8f"
$ws.Range("D2").Value = "obnoxious-serval@example.com"
$ws.Range("E2").Value = "Outrageous-Teira-Xerothermic-Iguana-2"
